# Simplifying databases and fixing input locator
# Target: SUPPLY sheet (type_dhw / type_el columns) + active-tab bookkeeping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUPPLY")

# --- type_dhw (column G): rows that held "T3" move to "T19", except the
#     HOTEL / RETAIL / FOODSTORE blocks, which move to "T20" instead.
#     (Rows left at "T0" are untouched; the LAB block at 182:193 is formula
#     driven off these same rows and recalculates on its own.)
$ws.Range("G24:G25").Value = "T19"
$ws.Range("G27:G37").Value = "T20"
$ws.Range("G39:G49").Value = "T19"
$ws.Range("G51:G61").Value = "T20"
$ws.Range("G63:G73").Value = "T20"
$ws.Range("G75:G97").Value = "T19"
$ws.Range("G99:G109").Value = "T19"
$ws.Range("G111:G121").Value = "T19"
$ws.Range("G123:G133").Value = "T19"
$ws.Range("G146:G157").Value = "T19"
$ws.Range("G170:G181").Value = "T19"
$ws.Range("G195:G205").Value = "T19"
$ws.Range("G207:G217").Value = "T19"
$ws.Range("G219:G229").Value = "T19"

# --- type_el (column H): every data row (2:229) moves from "T1" to "T24" ---
$ws.Range("H2:H229").Value = "T24"

# --- View bookkeeping: SUPPLY becomes the active/selected tab, replacing
#     INDOOR_COMFORT, with K224 as the last selected cell. ---
$ws.Activate()
$ws.Range("K224").Select()
